# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list on Sat Jun  3 19:29:15 UTC 2023 with GitHub Actions".
# D (Price) and E (Volume(1h)) columns hold scraped display text, not
# numbers, so D-column writes use a leading apostrophe to force Excel to
# keep them as literal text (preserving things like trailing zeros and
# the site's dotted thousands-separator prices such as 27.205.38).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.205.38"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "'1.894.13"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "'307.14"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").Value = "'0.5219"
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("D8").Value = "'0.3751"
$ws.Range("E8").Value = "  -0.78%  "
$ws.Range("D9").Value = "'0.07264"
$ws.Range("E9").Value = "  -0.51%  "
$ws.Range("D10").Value = "'21.17"
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("D11").Value = "'0.8982"
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("D12").Value = "'0.08152"
$ws.Range("E12").Value = "  +6.21%  "
$ws.Range("D13").Value = "'96.73"
$ws.Range("E13").Value = "  +2.05%  "
$ws.Range("D14").Value = "'1.891.69"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").Value = "'5.271"
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("D16").Value = "'1.004"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").Value = "'0.000008596"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("D18").Value = "'14.53"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("D20").Value = "'27.233.03"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "'5.080"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").Value = "'10.68"
$ws.Range("E22").Value = "  +0.45%  "
$ws.Range("D23").Value = "'6.397"
$ws.Range("E23").Value = "  -0.71%  "
$ws.Range("D24").Value = "'147.51"
$ws.Range("E24").Value = "  +1.14%  "
$ws.Range("D25").Value = "'2.287"
$ws.Range("E25").Value = "  -1.34%  "
$ws.Range("D26").Value = "'1.745"
$ws.Range("E26").Value = "  +0.97%  "
$ws.Range("D27").Value = "'18.19"
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("D28").Value = "'114.98"
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("D29").Value = "'4.905"
$ws.Range("E29").Value = "  -1.13%  "
$ws.Range("D30").Value = "'4.798"
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("D31").Value = "'0.09227"
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("D32").Value = "'0.05037"
$ws.Range("E32").Value = "  -0.73%  "
$ws.Range("D33").Value = "'0.7926"
$ws.Range("E33").Value = "  +1.45%  "
$ws.Range("D34").Value = "'1.215"
$ws.Range("E34").Value = "  -2.51%  "
$ws.Range("D35").Value = "'3.435"
$ws.Range("E35").Value = "  +3.89%  "
$ws.Range("D36").Value = "'2.950"
$ws.Range("E36").Value = "  -1.81%  "
$ws.Range("E37").Value = "  -1.61%  "
$ws.Range("D38").Value = "'0.5645"
$ws.Range("E38").Value = "  -0.70%  "
$ws.Range("D39").Value = "'0.01981"
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("D40").Value = "'1.074"
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D41").Value = "'8.921"
$ws.Range("E41").Value = "  -1.21%  "
$ws.Range("D42").Value = "'6.527"
$ws.Range("E42").Value = "  -1.58%  "
$ws.Range("D43").Value = "'115.05"
$ws.Range("E43").Value = "  -3.02%  "
$ws.Range("D44").Value = "'0.1516"
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("D45").Value = "'0.4866"
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("D47").Value = "'10.09"
$ws.Range("E47").Value = "  -1.58%  "
$ws.Range("D48").Value = "'1.615"
$ws.Range("E48").Value = "  +0.60%  "
$ws.Range("D49").Value = "'38.10"
$ws.Range("E49").Value = "  +1.82%  "
$ws.Range("D50").Value = "'63.27"
$ws.Range("E50").Value = "  -1.50%  "
$ws.Range("D51").Value = "'0.05944"
$ws.Range("E51").Value = "  +0.24%  "
